# Auto-generated edit script: applies the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "27.429.22"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.827.75"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  -0.98%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.38"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9634"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").Value = "1.826.05"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.871"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.072"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001024"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "27.430.22"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.320"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "2.064.88"
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.065"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.292"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09296"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9364"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.62%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.231"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.324"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05917"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02176"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.145"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5768"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1821"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.956"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.66%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.265"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5430"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.870"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06581"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("B51").Value = "PaxosStandard"
$ws.Range("C51").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -24.17%  "
